# Fix cells that were mistakenly left as text/shared-string placeholders
# (null/invalid time input) in column A/B/C/D/E/F of the "Data" sheet.
# Each affected cell is repaired by copying in the numeric date/time value
# of its adjacent cell in the same row, so both cells hold the same value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Each entry: cell to fix -> cell to copy the value from (same row)
$fixes = @(
    @{ Target = "A2";  Source = "B2"  },
    @{ Target = "F2";  Source = "E2"  },
    @{ Target = "F9";  Source = "E9"  },
    @{ Target = "F10"; Source = "E10" },
    @{ Target = "A13"; Source = "B13" },
    @{ Target = "B16"; Source = "A16" },
    @{ Target = "C22"; Source = "B22" },
    @{ Target = "F33"; Source = "E33" },
    @{ Target = "A41"; Source = "B41" },
    @{ Target = "D44"; Source = "E44" },
    @{ Target = "E48"; Source = "F48" },
    @{ Target = "E60"; Source = "F60" }
)

foreach ($fix in $fixes) {
    $value = $ws.Range($fix.Source).Value2
    $ws.Range($fix.Target).Value2 = $value
}
